$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting so numeric-looking
# strings (e.g. "1.00", thousands-dotted prices, percentages) are not
# auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.344.51"
$ws.Range("E2").Value = "  +5.54%  "
$ws.Range("D3").Value = "3.385.28"
$ws.Range("E3").Value = "  +6.05%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "576.77"
$ws.Range("E5").Value = "  +7.58%  "
$ws.Range("D6").Value = "154.24"
$ws.Range("E6").Value = "  +6.39%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.387.04"
$ws.Range("E8").Value = "  +5.87%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "7.46"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +6.50%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "3.975.79"
$ws.Range("E13").Value = "  +6.22%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  +6.56%  "
$ws.Range("D16").Value = "26.99"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Value = "63.453.37"
$ws.Range("E17").Value = "  +5.68%  "
$ws.Range("D18").Value = "3.384.21"
$ws.Range("E18").Value = "  +6.05%  "
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  +5.32%  "
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").Value = "388.55"
$ws.Range("E22").Value = "  +5.09%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "0.534"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "70.75"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "9.49"
$ws.Range("E26").Value = "  +10.82%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.179"
$ws.Range("E27").Value = "  +5.95%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("E28").Value = "  +17.97%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").Value = "  +7.64%  "
$ws.Range("D31").Value = "6.45"
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").Value = "23.08"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  +10.14%  "
$ws.Range("D34").Value = "5.57"
$ws.Range("E34").Value = "  +5.18%  "
$ws.Range("D35").Value = "6.72"
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("D36").Value = "1.48"
$ws.Range("E36").Value = "  +9.16%  "
$ws.Range("D37").Value = "158.60"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "27.68"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  +12.22%  "
$ws.Range("D40").Value = "0.0749"
$ws.Range("D41").Value = "2.876.48"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("D43").Value = "0.761"
$ws.Range("E43").Value = "  +5.75%  "
$ws.Range("D44").Value = "41.12"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "4.28"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  +7.88%  "
$ws.Range("D47").Value = "3.438.26"
$ws.Range("E47").Value = "  +6.28%  "
$ws.Range("D48").Value = "22.08"
$ws.Range("E48").Value = "  +6.97%  "
$ws.Range("D49").Value = "298.43"
$ws.Range("E49").Value = "  +12.48%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "6.33"
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.103"
$ws.Range("E51").Value = "  -0.99%  "
